$d = $word.ActiveDocument

$d.Content.Find.Execute("355×8=2840", $true, $false, $false, $false, $false, $true, 1, $false, "771×6=4626", 2) | Out-Null
$d.Content.Find.Execute("770×5=3850", $true, $false, $false, $false, $false, $true, 1, $false, "465×6=2790", 2) | Out-Null
$d.Content.Find.Execute("453×4=1812", $true, $false, $false, $false, $false, $true, 1, $false, "715×7=5005", 2) | Out-Null
$d.Content.Find.Execute("586×7=4102", $true, $false, $false, $false, $false, $true, 1, $false, "408×7=2856", 2) | Out-Null
$d.Content.Find.Execute("922×4=3688", $true, $false, $false, $false, $false, $true, 1, $false, "932×8=7456", 2) | Out-Null
$d.Content.Find.Execute("153×4=612", $true, $false, $false, $false, $false, $true, 1, $false, "744×2=1488", 2) | Out-Null
$d.Content.Find.Execute("763×8=6104", $true, $false, $false, $false, $false, $true, 1, $false, "622×9=5598", 2) | Out-Null
$d.Content.Find.Execute("272×3=816", $true, $false, $false, $false, $false, $true, 1, $false, "701×7=4907", 2) | Out-Null
$d.Content.Find.Execute("547×8=4376", $true, $false, $false, $false, $false, $true, 1, $false, "156×3=468", 2) | Out-Null
$d.Content.Find.Execute("152×6=912", $true, $false, $false, $false, $false, $true, 1, $false, "693×8=5544", 2) | Out-Null
$d.Content.Find.Execute("881×2=1762", $true, $false, $false, $false, $false, $true, 1, $false, "518×2=1036", 2) | Out-Null
$d.Content.Find.Execute("210×8=1680", $true, $false, $false, $false, $false, $true, 1, $false, "424×9=3816", 2) | Out-Null
$d.Content.Find.Execute("774×3=2322", $true, $false, $false, $false, $false, $true, 1, $false, "501×7=3507", 2) | Out-Null
$d.Content.Find.Execute("190×5=950", $true, $false, $false, $false, $false, $true, 1, $false, "700×2=1400", 2) | Out-Null
$d.Content.Find.Execute("614×8=4912", $true, $false, $false, $false, $false, $true, 1, $false, "763×4=3052", 2) | Out-Null
$d.Content.Find.Execute("913×9=8217", $true, $false, $false, $false, $false, $true, 1, $false, "754×2=1508", 2) | Out-Null
$d.Content.Find.Execute("232×8=1856", $true, $false, $false, $false, $false, $true, 1, $false, "807×7=5649", 2) | Out-Null
$d.Content.Find.Execute("663×9=5967", $true, $false, $false, $false, $false, $true, 1, $false, "317×9=2853", 2) | Out-Null
$d.Content.Find.Execute("722×2=1444", $true, $false, $false, $false, $false, $true, 1, $false, "112×2=224", 2) | Out-Null
$d.Content.Find.Execute("604×4=2416", $true, $false, $false, $false, $false, $true, 1, $false, "421×9=3789", 2) | Out-Null
$d.Content.Find.Execute("649×8=5192", $true, $false, $false, $false, $false, $true, 1, $false, "921×6=5526", 2) | Out-Null
$d.Content.Find.Execute("899×7=6293", $true, $false, $false, $false, $false, $true, 1, $false, "514×3=1542", 2) | Out-Null
$d.Content.Find.Execute("902×5=4510", $true, $false, $false, $false, $false, $true, 1, $false, "591×2=1182", 2) | Out-Null
$d.Content.Find.Execute("320×4=1280", $true, $false, $false, $false, $false, $true, 1, $false, "322×5=1610", 2) | Out-Null
$d.Content.Find.Execute("450×4=1800", $true, $false, $false, $false, $false, $true, 1, $false, "914×8=7312", 2) | Out-Null
